$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the duplicated bold title paragraph near the end of the
#    document, and rewrite the following italic paragraph's text with
#    the new image-generation prompt. Do this first, while paragraph
#    indices still reflect the pristine document layout.
# ------------------------------------------------------------------
$dupTitlePara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if (($txt -eq "Play Alice & The Mad Tea Party Slot for Free | Review") -and ($p.Range.Font.Bold)) {
        $dupTitlePara = $p
    }
}

if ($dupTitlePara -ne $null) {
    $idx = $dupTitlePara.Index
    $nextPara = $d.Paragraphs($idx + 1)
    $delStart = $dupTitlePara.Range.Start
    $delEnd = $nextPara.Range.Start
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$oldDescription = "Read the review of Alice & The Mad Tea Party online slot game. Play for free with interesting winning potentials and special features."
$newDescription = "Prompt: Create a feature image for Alice & The Mad Tea Party that showcases the fun and whimsical vibe of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be surrounded by the characters from the game, including Alice, the Mad Hatter, and the White Rabbit. The background should be a colorful and playful representation of Wonderland, with colorful mushrooms, teacups, and flowers. The overall tone of the image should be upbeat and lighthearted to reflect the exciting gameplay and potential for big wins."

$lastPara.Range.Find.Execute($oldDescription, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)

# ------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the H1
#    title at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$insertionPoint = $d.Range($metaStart, $metaStart)
$insertionPoint.InsertAfter("Meta description: Read the review of Alice & The Mad Tea Party online slot game. Play for free with interesting winning potentials and special features.")

$boldLen = ("Meta description").Length
$boldRange = $d.Range($metaStart, $metaStart + $boldLen)
$boldRange.Font.Bold = 1

Write-Output "done"
